$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.122.90"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").Value = "2.328.87"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.73"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.62%  "
$ws.Range("E6").Value = "  +0.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.507"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.24%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.504"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.62"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.42"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +7.29%  "
$ws.Range("E12").Value = "  +1.19%  "
$ws.Range("E13").Value = "  +0.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.93"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.85%  "
$ws.Range("D15").Value = "2.692.58"
$ws.Range("E15").Value = "  +1.10%  "
$ws.Range("D16").Value = "2.323.37"
$ws.Range("E16").Value = "  +0.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.788"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "43.027.27"
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.54"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.03%  "
$ws.Range("D20").Value = "0.0₃0901"
$ws.Range("E20").Value = "  -0.30%  "
$ws.Range("E21").Value = "  +0.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.00"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.55"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.20"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.25%  "
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.45"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.78%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.91"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "166.41"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.39%  "
$ws.Range("E29").Value = "  +1.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.14"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.19"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.01"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.74%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.55"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -8.42%  "
$ws.Range("E36").Value = "  -1.27%  "
$ws.Range("E37").Value = "  +1.39%  "
$ws.Range("E38").Value = "  +0.18%  "
$ws.Range("E39").Value = "  +2.37%  "
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("E41").Value = "  -0.45%  "
$ws.Range("D42").Value = "1.998.34"
$ws.Range("E42").Value = "  -0.63%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.74"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.87%  "
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.21"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.56%  "
$ws.Range("E46").Value = "  -2.66%  "
$ws.Range("E47").Value = "  -0.77%  "
$ws.Range("D48").Value = "2.558.81"
$ws.Range("E48").Value = "  +1.05%  "
$ws.Range("E49").Value = "  -0.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.78"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.99"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.46%  "
